$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.416.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = "'1.723.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'243.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = "'0.4937"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").Value = "'0.2615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.44%  '
$ws.Range("D9").Value = "'0.06214"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").Value = "'1.727.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").Value = "'0.06991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").Value = "'15.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.21%  '
$ws.Range("D13").Value = "'4.554"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").Value = "'0.6003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").Value = "'77.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = "'26.414.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").Value = "'0.000007203"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("D20").Value = "'11.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = "'1.944.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("D22").Value = "'4.474"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("D23").Value = "'8.586"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'5.165"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.43%  '
$ws.Range("D25").Value = "'137.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").Value = "'1.399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("D28").Value = "'107.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("D29").Value = "'1.722"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").Value = "'0.08015"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").Value = "'0.04522"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = "'0.9993"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("D36").Value = "'0.9994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("D37").Value = "'0.6262"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("D38").Value = "'0.9444"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.55%  '
$ws.Range("D39").Value = "'2.391"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").Value = "'1.947"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.93%  '
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = "'0.01485"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = "'99.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.25%  '
$ws.Range("D44").Value = "'5.311"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.83%  '
$ws.Range("D45").Value = "'0.3857"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").Value = "'6.812"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.42%  '
$ws.Range("D47").Value = "'0.1170"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = "'0.05369"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '
$ws.Range("D49").Value = "'7.786"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("D51").Value = "'1.233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.58%  '
